# Apply text updates described in the diff:
#  - Slide 2: update Bulk/Single-cell pipeline summary labels (v2 pipeline overview box)
#  - Slide 8: update single-cell pipeline title

$p = $ppt.ActivePresentation

# --- Slide 2 ("Pipeline Overview") ---
$s2 = $p.Slides.Item(2)

# Shape "Text 17": "Bulk: 2-Step Process" -> "Bulk RNA-seq: 6-Agent Pipeline"
$s2.Shapes.Item(18).TextFrame.TextRange.Text = "Bulk RNA-seq: 6-Agent Pipeline"

# Shape "Text 18": "Step 1: DEG Analysis, Step 2: Interpretation" -> "DEG -> Network -> Pathway -> Validation -> Viz -> Report"
$s2.Shapes.Item(19).TextFrame.TextRange.Text = "DEG → Network → Pathway → Validation → Viz → Report"

# Shape "Text 20": "Single-cell: 1-Step Process" -> "Single-cell: Scanpy Pipeline"
$s2.Shapes.Item(21).TextFrame.TextRange.Text = "Single-cell: Scanpy Pipeline"

# Shape "Text 21": "QC - Clustering - Annotation - Report" -> "QC -> Clustering -> Annotation -> Markers -> Report"
$s2.Shapes.Item(22).TextFrame.TextRange.Text = "QC → Clustering → Annotation → Markers → Report"

# --- Slide 8 ("Single-cell Pipeline") ---
$s8 = $p.Slides.Item(8)

# Shape "Text 1": "Single-cell Pipeline (1-Step)" -> "Single-cell RNA-seq Pipeline"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "Single-cell RNA-seq Pipeline"
